$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 6797  # was 6794
$ws.Range("G2").Value = 80  # was 70
$ws.Range("F3").Value = 819  # was 818
$ws.Range("G3").Value = 68  # was 58
$ws.Range("F5").Value = 144  # was 142
$ws.Range("F6").Value = 11  # was 10
$ws.Range("F11").Value = 25  # was 24
$ws.Range("F15").Value = 704  # was 703
$ws.Range("F19").Value = 120  # was 119
$ws.Range("F20").Value = 539  # was 536
$ws.Range("F22").Value = 567  # was 563
$ws.Range("F23").Value = 12  # was 11
$ws.Range("F26").Value = 1066  # was 1065
$ws.Range("F27").Value = 1495  # was 1497
$ws.Range("F29").Value = 543  # was 540
$ws.Range("F30").Value = 463  # was 462
$ws.Range("F32").Value = 94  # was 93
$ws.Range("F36").Value = 2365  # was 2363
$ws.Range("F38").Value = 1261  # was 1260
$ws.Range("F41").Value = 3873  # was 3872

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 749  # was 748
$ws.Range("F7").Value = 1035  # was 1034
$ws.Range("F17").Value = 4129  # was 4127
$ws.Range("F20").Value = 25  # was 24
$ws.Range("F23").Value = 238  # was 237
$ws.Range("F31").Value = 1705  # was 1704

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 1256  # was 1255
$ws.Range("F6").Value = 446  # was 445
$ws.Range("F8").Value = 970  # was 969

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1256  # was 1255
$ws.Range("F5").Value = 446  # was 445
$ws.Range("F7").Value = 970  # was 969
$ws.Range("F8").Value = 6797  # was 6794
$ws.Range("G8").Value = 80  # was 70
$ws.Range("F10").Value = 819  # was 818
$ws.Range("G10").Value = 68  # was 58
$ws.Range("F11").Value = 749  # was 748
$ws.Range("F12").Value = 144  # was 142
$ws.Range("F13").Value = 11  # was 10
$ws.Range("F17").Value = 25  # was 24
$ws.Range("F21").Value = 704  # was 703
$ws.Range("F28").Value = 120  # was 119
$ws.Range("F29").Value = 539  # was 536
$ws.Range("F30").Value = 567  # was 563
$ws.Range("F34").Value = 1066  # was 1065
$ws.Range("F35").Value = 1495  # was 1497
$ws.Range("F37").Value = 543  # was 540
$ws.Range("F38").Value = 463  # was 462
$ws.Range("F45").Value = 2365  # was 2363
$ws.Range("F47").Value = 1705  # was 1704
$ws.Range("F48").Value = 1705  # was 1704
$ws.Range("F49").Value = 1261  # was 1260
$ws.Range("F51").Value = 3873  # was 3872
